# "Se añaden planos cable datos p10 y p6.67"
# Bump two item quantities and retarget the sheet's print setup for the
# newly-added cable-plan drawings: tighter margins, 59% print scale
# (fit-to-page), Page Break Preview at 60% zoom, scrolled down with the
# selection parked on A8.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Cell value updates ---
$ws.Range("B10").Value = 10
$ws.Range("B11").Value = 10

# --- Page setup: margins, print scale, fit-to-page ---
$ps = $ws.PageSetup

$ps.LeftMargin   = $excel.InchesToPoints(0.25)
$ps.RightMargin  = $excel.InchesToPoints(0.25)
$ps.TopMargin    = $excel.InchesToPoints(0.75)
$ps.BottomMargin = $excel.InchesToPoints(0.75)
$ps.HeaderMargin = $excel.InchesToPoints(0.3)
$ps.FooterMargin = $excel.InchesToPoints(0.3)

# Record the manual scale, then engage "fit to page" (1 page wide x 1 tall).
$ps.Zoom = 59
$ps.FitToPagesWide = 1
$ps.FitToPagesTall = 1

# --- Window/view: Page Break Preview, zoomed to 60%, scrolled to row 73 ---
$win = $excel.ActiveWindow
$win.View = 2
$win.Zoom = 60
$win.ScrollRow = 73
$win.ScrollColumn = 1

# --- Selection moves to A8 ---
$ws.Range("A8").Select() | Out-Null
